$wb = $excel.ActiveWorkbook

# Target sheet: "Đơn phụ phẫu 2" (6th worksheet) — the diff adds a new row 2
# (a blank data row below the header row) and extends the sheet dimension
# from A1:T1 to A1:T2.
$ws = $wb.Worksheets.Item(6)

# Columns that hold numeric/currency values end up with an explicit 0 on
# the new blank row: I (Đơn giá gốc), K (Upsale), L (Đơn giá),
# M (Thanh toán lần đầu), N (Trả sau), O (Đã thanh toán), P (Dư nợ).
$numericCols = @(9, 11, 12, 13, 14, 15, 16)
foreach ($col in $numericCols) {
    $ws.Cells.Item(2, $col).Value = 0
}

# Remaining columns (A, B, C, D, E, F, G, H, J, Q, R, S, T) stay blank on
# the new row, but the row/cells themselves still need to exist. Touching
# a harmless, already-default formatting property materializes the cell
# in the sheet without giving it a value or changing any style.
$blankCols = @(1, 2, 3, 4, 5, 6, 7, 8, 10, 17, 18, 19, 20)
foreach ($col in $blankCols) {
    $ws.Cells.Item(2, $col).Font.Bold = $false
}
